$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-44 is updated from 2025-03-17 (45733)
# to 2025-03-18 (45734) for every data row.
for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 3).Value = 45734
}
